$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C column) date for rows 2 through 11 from 2023-10-25 (45224) to 2023-11-03 (45233)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
